$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P mirrors column O's formatting; copy the styles first, then
# overwrite the values with the new 2022 figures.
$ws.Range("O3:O5").Copy()
$ws.Range("P3:P5").PasteSpecial(-4122)

$ws.Range("P3").Value = 2022
$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 2130.4

# Match the author's final selection.
$ws.Range("P6").Select()
